# Updated cryptos list (Price/Volume(1h) refresh + one ONDO/OKB row-order swap).
# All target cells are plain text (t="inlineStr"/shared-string Text cells),
# so every write goes through Set-TextCell to stop Excel from re-typing
# number-looking strings (e.g. "1.00", "41.63") as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces Excel's text-literal parsing so values
    # like "1.00" or "41.63" stay text instead of becoming numbers;
    # resetting the style afterwards drops the quote-prefix flag again
    # so the cell keeps its original (unstyled) appearance.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "64.106.22"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "3.519.87"
$ws.Range("E4").Value = "  -0.01%  "
Set-TextCell "D5" "586.69"
$ws.Range("E5").Value = "  +0.12%  "
Set-TextCell "D6" "134.14"
$ws.Range("D7").Value = "3.519.82"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "4.123.11"
$ws.Range("E13").Value = "  +0.23%  "
Set-TextCell "D14" "27.44"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "3.519.52"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "64.158.18"
$ws.Range("E18").Value = "  -1.26%  "
Set-TextCell "D19" "9.78"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("E20").Value = "  -2.86%  "
$ws.Range("E21").Value = "  -0.64%  "
Set-TextCell "D22" "382.97"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").Value = "3.664.03"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -1.46%  "
Set-TextCell "D28" "0.0000114"
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  -1.04%  "
Set-TextCell "D31" "1.00"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "3.535.76"
$ws.Range("E34").Value = "  +0.36%  "
Set-TextCell "D36" "23.56"
$ws.Range("E36").Value = "  -1.99%  "
Set-TextCell "D37" "0.146"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("E39").Value = "  +0.31%  "
Set-TextCell "D40" "1.55"
$ws.Range("E40").Value = "  -0.62%  "
Set-TextCell "D41" "158.40"
$ws.Range("E41").Value = "  -6.47%  "
Set-TextCell "D42" "0.0788"
$ws.Range("E42").Value = "  -2.11%  "
Set-TextCell "D43" "26.61"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D46" "41.63"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell "D47" "1.21"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").Value = "2.482.64"
$ws.Range("E50").Value = "  +0.47%  "
Set-TextCell "D51" "6.79"
$ws.Range("E51").Value = "  -0.81%  "
